$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply wrap-text formatting to the new row's columns first (creates the
# plain "wrapText only" style used by columns A/B), then bump the font size
# on C:E (creates the "font 8 + wrapText" style used by columns C/D/E).
$ws.Range("A2:B2").WrapText = $true
$ws.Range("C2:E2").WrapText = $true
$ws.Range("C2:E2").Font.Size = 8

# Fill in the new row of data (set in the order that matches the target
# shared-string table: English line, filename, number, Russian line,
# "converted" line).
$ws.Range("C2").Value = " I wish you the best of fortunes\non your expedition. Hee-hee!"
$ws.Range("A2").Value = "SCRIPT/T01P01A/enter05.ssb"
$ws.Range("B2").Value = 336
$ws.Range("D2").Value = " Желаю вам всего самого\nнаилучшего в экспедиции. Хи-хи!"
$ws.Range("E2").Value = " Çåìàý âàí âòåãï òàíïãï\nîàéìôœšåãï â üëòðåäéøéé. Öé-öé!"

# Row height for the new row.
$ws.Rows.Item(2).RowHeight = 43.2

# Update the active selection to match the target workbook state.
$ws.Range("C1").Select()
